$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pop times")

# Fix "quantifer" -> "quantifier" typo in existing cells (text unchanged otherwise)
$ws.Range("D4").Value = "Only in years between campaigns, varying the CD quantifier from population x 0% to population x 40%"
$ws.Range("D3").Value = "2022-2032 varying the CD quantifier from population x 0% to population x 50%"

# Rename header label
$ws.Range("E1").Value = "Number of different iterations per scenario"

# Update the iteration descriptions and counts
$ws.Range("B5").Value = "In 2022, 2025, 2028, 2031, 2034, varying from population / 0.1-2.0"
$ws.Range("E5").Value = 21

$ws.Range("B6").Value = "In 2022, 2024, 2026, 2028, 2030, 2032, 2034 varying from population / 0.5-2.0"
$ws.Range("E6").Value = 16

# Update active selection
$ws.Range("G5").Select()
